# ---------------------------------------------------------------------------
# Applies two changes to the presentation:
#
#  1. Slide 5's table (graphicFrame "Google Shape;122;p17") switches its
#     table style ("tableStyleId") from the deck's custom style
#     {5AF98349-6BCA-41DD-90F6-C31DFD8827F7} to the built-in style
#     {0889060E-1D78-4553-BF7E-A0203FAA7392}.
#
#  2. The presentation's theme colour scheme (the one used by the slide
#     master, backing ppt/theme/theme1.xml) is switched from the
#     "Integral"/"Red Violet" palette to the standard Office theme palette
#     (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -----------------------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{0889060E-1D78-4553-BF7E-A0203FAA7392}")
    }
}

# --- 2. Theme colours -----------------------------------------------------
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$scheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $scheme.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
